# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.451.91"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "'1.938.91"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'243.53"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'57.45"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.361"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "'0.0853"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "'2.220.91"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "'21.41"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "'0.812"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "'13.39"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "'5.16"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "'1.930.28"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "'36.370.28"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'69.15"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "'0.0₃0862"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "'227.36"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").Value = "'4.97"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'2.35"
$ws.Range("E24").Value = "  -5.65%  "
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'9.18"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").Value = "'160.61"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").Value = "'0.132"
$ws.Range("E28").Value = "  +14.88%  "
$ws.Range("D29").Value = "'19.16"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("D32").Value = "'4.58"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").Value = "'0.0615"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "'4.17"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'6.19"
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").Value = "'2.18"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").Value = "'3.11"
$ws.Range("E39").Value = "  +8.74%  "
$ws.Range("D40").Value = "'0.0980"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").Value = "'2.95"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").Value = "'15.82"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.03"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'1.333.14"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'7.15"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "'85.80"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'2.112.32"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'43.14"
$ws.Range("E51").Value = "  -2.42%  "
